$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"
$t.Cell(4, 1).Range.Text = "15391"

$t.Cell(6, 1).Range.Text = "0.27930"
$t.Cell(7, 1).Range.Text = "0.06101"
$t.Cell(8, 1).Range.Text = "0.00762"
$t.Cell(9, 1).Range.Text = "0.24525"
$t.Cell(10, 1).Range.Text = "0.24525"
$t.Cell(11, 1).Range.Text = "0.27930"
$t.Cell(12, 1).Range.Text = "46.85050"

$t.Cell(44, 1).Range.Text = "98.52"
$t.Cell(45, 1).Range.Text = "46.85"
$t.Cell(46, 1).Range.Text = "3172"
